{"js": "// The template text \"Application number: {{ proposal.lodgement_number }}\"\n// needs \"proposal\" replaced with \"application\" inside the merge field\n// expression, so the final text reads:\n// \"Application number: {{ application.lodgement_number }}\"\nconst results = context.document.body.search(\"proposal.lodgement_number\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"application.lodgement_number\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace \"proposal\" with \"application\" in the merge-field expression so\n# \"Application number: {{ proposal.lodgement_number }}\" becomes\n# \"Application number: {{ application.lodgement_number }}\".\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"proposal.lodgement_number\"\n$find.Replacement.Text = \"application.lodgement_number\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n"}
